# Auto-generated Excel COM-interop script to apply market-price/profit updates
# across multiple worksheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR) per the diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 25970
$ws.Range("J57").Value = 25970
$ws.Range("L57").Value = 77910
$ws.Range("N57").Value = -78908
$ws.Range("H70").Value = 1452.9412
$ws.Range("I70").Value = 877.7778
$ws.Range("J70").Value = 2100
$ws.Range("K70").Value = 2633.3334
$ws.Range("L70").Value = 6300
$ws.Range("M70").Value = -2363.3334
$ws.Range("N70").Value = -6840
$ws.Range("H73").Value = 1452.9412
$ws.Range("I73").Value = 877.7778
$ws.Range("J73").Value = 2100
$ws.Range("K73").Value = 2633.3334
$ws.Range("L73").Value = 6300
$ws.Range("M73").Value = -1697.3334
$ws.Range("N73").Value = -8172
$ws.Range("H86").Value = 76926770
$ws.Range("I86").Value = 5401.5
$ws.Range("J86").Value = 90912470
$ws.Range("K86").Value = 5401.5
$ws.Range("L86").Value = 90912470
$ws.Range("M86").Value = -4278.5
$ws.Range("N86").Value = -90914716
$ws.Range("H88").Value = 1765486.5
$ws.Range("I88").Value = 3867.3333
$ws.Range("J88").Value = 2646296.2
$ws.Range("K88").Value = 3867.3333
$ws.Range("L88").Value = 2646296.2
$ws.Range("M88").Value = -3461.3333
$ws.Range("N88").Value = -2647108.2
$ws.Range("H89").Value = 76926770
$ws.Range("I89").Value = 5401.5
$ws.Range("J89").Value = 90912470
$ws.Range("K89").Value = 27007.5
$ws.Range("L89").Value = 454562350
$ws.Range("M89").Value = -21391.5
$ws.Range("N89").Value = -454573582
$ws.Range("H91").Value = 1765486.5
$ws.Range("I91").Value = 3867.3333
$ws.Range("J91").Value = 2646296.2
$ws.Range("K91").Value = 3867.3333
$ws.Range("L91").Value = 2646296.2
$ws.Range("M91").Value = -2463.3333
$ws.Range("N91").Value = -2649104.2
$ws.Range("H99").Value = 433
$ws.Range("I99").Value = 200
$ws.Range("J99").Value = 666
$ws.Range("K99").Value = 600
$ws.Range("L99").Value = 1998
$ws.Range("M99").Value = 898
$ws.Range("N99").Value = -4994
$ws.Range("H137").Value = 14707438
$ws.Range("I137").Value = 1247.0526
$ws.Range("J137").Value = 33335280
$ws.Range("K137").Value = 3741.1578
$ws.Range("L137").Value = 100005840
$ws.Range("M137").Value = -1191.1578
$ws.Range("N137").Value = -100010940

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 966.8095
$ws.Range("I61").Value = 788.41174
$ws.Range("J61").Value = 1725
$ws.Range("K61").Value = 788.41174
$ws.Range("L61").Value = 1725
$ws.Range("M61").Value = -576.41174
$ws.Range("N61").Value = -2149
$ws.Range("H97").Value = 582.3077
$ws.Range("I97").Value = 406.08694
$ws.Range("J97").Value = 1933.3334
$ws.Range("K97").Value = 406.08694
$ws.Range("L97").Value = 1933.3334
$ws.Range("M97").Value = 89.91305999999997
$ws.Range("N97").Value = -2925.3334
$ws.Range("H102").Value = 2588.3333
$ws.Range("I102").Value = 2706
$ws.Range("K102").Value = 2706
$ws.Range("M102").Value = -1084
$ws.Range("H132").Value = 1655.8948
$ws.Range("I132").Value = 1215.9231
$ws.Range("J132").Value = 2609.1667
$ws.Range("K132").Value = 3647.7693
$ws.Range("L132").Value = 7827.500100000001
$ws.Range("M132").Value = -1117.7693
$ws.Range("N132").Value = -12887.5001
$ws.Range("H136").Value = 966.8095
$ws.Range("I136").Value = 788.41174
$ws.Range("J136").Value = 1725
$ws.Range("K136").Value = 2365.23522
$ws.Range("L136").Value = 5175
$ws.Range("M136").Value = 184.76478
$ws.Range("N136").Value = -10275

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1864.3704
$ws.Range("I132").Value = 1228.8
$ws.Range("K132").Value = 3686.4
$ws.Range("M132").Value = -1156.4
$ws.Range("H133").Value = 48950
$ws.Range("J133").Value = 48950
$ws.Range("L133").Value = 48950
$ws.Range("N133").Value = -54010
$ws.Range("H139").Value = 41425
$ws.Range("I139").Value = 26000
$ws.Range("J139").Value = 46566.668
$ws.Range("K139").Value = 26000
$ws.Range("L139").Value = 46566.668
$ws.Range("M139").Value = -20860
$ws.Range("N139").Value = -56846.668
$ws.Range("H140").Value = 49100
$ws.Range("J140").Value = 49100
$ws.Range("L140").Value = 49100
$ws.Range("N140").Value = -59460

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 1842.25
$ws.Range("K42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("M42").Value = 5526.75
$ws.Range("N42").Value = -6594.75
$ws.Range("H64").Value = 2072.4
$ws.Range("I64").Value = 1008
$ws.Range("J64").Value = 2528.5715
$ws.Range("K64").Value = 3024
$ws.Range("L64").Value = 7585.7145
$ws.Range("M64").Value = -2754
$ws.Range("N64").Value = -8125.7145
$ws.Range("H67").Value = 2072.4
$ws.Range("I67").Value = 1008
$ws.Range("J67").Value = 2528.5715
$ws.Range("K67").Value = 3024
$ws.Range("L67").Value = 7585.7145
$ws.Range("M67").Value = -2088
$ws.Range("N67").Value = -9457.7145
$ws.Range("H113").Value = 533.4433
$ws.Range("I113").Value = 476.3158
$ws.Range("J113").Value = 547.35895
$ws.Range("K113").Value = 1428.9474
$ws.Range("L113").Value = 1642.07685
$ws.Range("M113").Value = 741.0526
$ws.Range("N113").Value = -5982.07685
$ws.Range("H129").Value = 2797.8
$ws.Range("I129").Value = 735.6
$ws.Range("J129").Value = 8984.4
$ws.Range("K129").Value = 2206.8
$ws.Range("L129").Value = 26953.2
$ws.Range("M129").Value = 2793.2
$ws.Range("N129").Value = -36953.2

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3242.8572
$ws.Range("I80").Value = 3300
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 3300
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -2302
$ws.Range("N80").Value = -4896
$ws.Range("H83").Value = 3242.8572
$ws.Range("I83").Value = 3300
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 16500
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -11508
$ws.Range("N83").Value = -24484
$ws.Range("H97").Value = 1246.5927
$ws.Range("I97").Value = 996.8
$ws.Range("J97").Value = 1960.2858
$ws.Range("K97").Value = 996.8
$ws.Range("L97").Value = 1960.2858
$ws.Range("M97").Value = -500.8
$ws.Range("N97").Value = -2952.2858
$ws.Range("H137").Value = 54803.332
$ws.Range("J137").Value = 54803.332
$ws.Range("L137").Value = 54803.332
$ws.Range("N137").Value = -65003.332

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2474.0908
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 2173.5715
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 2173.5715
$ws.Range("M82").Value = -2639
$ws.Range("N82").Value = -2895.5715
$ws.Range("H85").Value = 2474.0908
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 2173.5715
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 2173.5715
$ws.Range("M85").Value = -1752
$ws.Range("N85").Value = -4669.5715
$ws.Range("H100").Value = 13890688
$ws.Range("I100").Value = 37038370
$ws.Range("J100").Value = 2079.8
$ws.Range("K100").Value = 37038370
$ws.Range("L100").Value = 2079.8
$ws.Range("M100").Value = -37037829
$ws.Range("N100").Value = -3161.8
$ws.Range("H139").Value = 32400
$ws.Range("J139").Value = 32400
$ws.Range("L139").Value = 32400
$ws.Range("N139").Value = -42680

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1812.5714
$ws.Range("I81").Value = 910.8570999999999
$ws.Range("J81").Value = 2714.2856
$ws.Range("K81").Value = 1821.7142
$ws.Range("L81").Value = 5428.5712
$ws.Range("M81").Value = -760.7141999999999
$ws.Range("N81").Value = -7550.5712
$ws.Range("H84").Value = 1812.5714
$ws.Range("I84").Value = 910.8570999999999
$ws.Range("J84").Value = 2714.2856
$ws.Range("K84").Value = 9108.571
$ws.Range("L84").Value = 27142.856
$ws.Range("M84").Value = -3804.571
$ws.Range("N84").Value = -37750.856
$ws.Range("H139").Value = 49526.668
$ws.Range("J139").Value = 49526.668
$ws.Range("L139").Value = 49526.668
$ws.Range("N139").Value = -59806.668

Write-Output "Applied all Bahamut_Profits cell updates"